$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = '76.629.73'
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.Value = '  +1.06%  '
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.Value = '2.895.26'
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.Value = '  +7.74%  '
$c.Style = "Normal"

$c = $ws.Range("E4")
$c.Value = '  +0.05%  '
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.Value = '''196.36'
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.Value = '  +4.53%  '
$c.Style = "Normal"

$c = $ws.Range("D6")
$c.Value = '''599.29'
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.Value = '  +1.87%  '
$c.Style = "Normal"

$c = $ws.Range("E7")
$c.Value = '  +0.06%  '
$c.Style = "Normal"

$c = $ws.Range("D8")
$c.Value = '''0.556'
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.Value = '  +2.98%  '
$c.Style = "Normal"

$c = $ws.Range("D9")
$c.Value = '''0.192'
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.Value = '  -2.16%  '
$c.Style = "Normal"

$c = $ws.Range("D10")
$c.Value = '2.895.42'
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.Value = '  +7.84%  '
$c.Style = "Normal"

$c = $ws.Range("D11")
$c.Value = '''0.398'
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.Value = '  +10.76%  '
$c.Style = "Normal"

$c = $ws.Range("E12")
$c.Value = '  -1.84%  '
$c.Style = "Normal"

$c = $ws.Range("D13")
$c.Value = '''4.91'
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.Value = '  +3.93%  '
$c.Style = "Normal"

$c = $ws.Range("D14")
$c.Value = '3.421.39'
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.Value = '  +7.93%  '
$c.Style = "Normal"

$c = $ws.Range("D15")
$c.Value = '76.526.90'
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.Value = '  +1.26%  '
$c.Style = "Normal"

$c = $ws.Range("D16")
$c.Value = '''27.49'
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.Value = '  +3.46%  '
$c.Style = "Normal"

$c = $ws.Range("D17")
$c.Value = '''0.0000189'
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.Value = '  +0.16%  '
$c.Style = "Normal"

$c = $ws.Range("D18")
$c.Value = '2.893.32'
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.Value = '  +7.37%  '
$c.Style = "Normal"

$c = $ws.Range("D19")
$c.Value = '''9.00'
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.Value = '  -3.02%  '
$c.Style = "Normal"

$c = $ws.Range("D20")
$c.Value = '''12.58'
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.Value = '  +5.08%  '
$c.Style = "Normal"

$c = $ws.Range("D21")
$c.Value = '''382.82'
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.Value = '  +2.53%  '
$c.Style = "Normal"

$c = $ws.Range("D22")
$c.Value = '''2.32'
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.Value = '  +1.33%  '
$c.Style = "Normal"

$c = $ws.Range("E23")
$c.Value = '  +1.15%  '
$c.Style = "Normal"

$c = $ws.Range("D24")
$c.Value = '''71.87'
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.Value = '  +2.56%  '
$c.Style = "Normal"

$c = $ws.Range("E25")
$c.Value = '  +0.50%  '
$c.Style = "Normal"

$c = $ws.Range("D26")
$c.Value = '3.043.04'
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.Value = '  +7.58%  '
$c.Style = "Normal"

$c = $ws.Range("D27")
$c.Value = '''4.22'
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.Value = '  +0.66%  '
$c.Style = "Normal"

$c = $ws.Range("D28")
$c.Value = '''9.75'
$c.Style = "Normal"
$c = $ws.Range("E28")
$c.Value = '  +3.79%  '
$c.Style = "Normal"

$c = $ws.Range("D29")
$c.Value = '''0.0000105'
$c.Style = "Normal"
$c = $ws.Range("E29")
$c.Value = '  +10.57%  '
$c.Style = "Normal"

$c = $ws.Range("D30")
$c.Value = '''1.00'
$c.Style = "Normal"
$c = $ws.Range("E30")
$c.Value = '  -0.10%  '
$c.Style = "Normal"

$c = $ws.Range("E31")
$c.Value = '  -0.26%  '
$c.Style = "Normal"

$c = $ws.Range("D32")
$c.Value = '''510.53'
$c.Style = "Normal"

$c = $ws.Range("E33")
$c.Value = '  -0.01%  '
$c.Style = "Normal"

$c = $ws.Range("E34")
$c.Value = '  +3.20%  '
$c.Style = "Normal"

$c = $ws.Range("D35")
$c.Value = '''1.00'
$c.Style = "Normal"
$c = $ws.Range("E35")
$c.Value = '  +0.11%  '
$c.Style = "Normal"

$c = $ws.Range("D36")
$c.Value = '''167.91'
$c.Style = "Normal"
$c = $ws.Range("E36")
$c.Value = '  +2.88%  '
$c.Style = "Normal"

$c = $ws.Range("D37")
$c.Value = '''20.12'
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.Value = '  +4.46%  '
$c.Style = "Normal"

$c = $ws.Range("D38")
$c.Value = '''0.118'
$c.Style = "Normal"
$c = $ws.Range("E38")
$c.Value = '  -1.32%  '
$c.Style = "Normal"

$c = $ws.Range("D39")
$c.Value = '''19.53'
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.Value = '  +0.76%  '
$c.Style = "Normal"

$c = $ws.Range("D40")
$c.Value = '''183.68'
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.Value = '  +8.28%  '
$c.Style = "Normal"

$c = $ws.Range("E41")
$c.Value = '  -0.03%  '
$c.Style = "Normal"

$c = $ws.Range("D42")
$c.Value = '''0.346'
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.Value = '  +4.91%  '
$c.Style = "Normal"

$c = $ws.Range("D43")
$c.Value = '''5.09'
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.Value = '  +1.49%  '
$c.Style = "Normal"

$c = $ws.Range("D44")
$c.Value = '''1.68'
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.Value = '  -1.25%  '
$c.Style = "Normal"

$c = $ws.Range("D45")
$c.Value = '''0.0924'
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.Value = '  +9.32%  '
$c.Style = "Normal"

$c = $ws.Range("E46")
$c.Value = '  +3.34%  '
$c.Style = "Normal"

$c = $ws.Range("D47")
$c.Value = '''40.23'
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.Value = '  +2.84%  '
$c.Style = "Normal"

$c = $ws.Range("D48")
$c.Value = '''2.35'
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.Value = '  -1.18%  '
$c.Style = "Normal"

$c = $ws.Range("D49")
$c.Value = '''0.698'
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.Value = '  +17.97%  '
$c.Style = "Normal"

$c = $ws.Range("D50")
$c.Value = '''0.581'
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.Value = '  +8.31%  '
$c.Style = "Normal"

$c = $ws.Range("D51")
$c.Value = '''3.77'
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.Value = '  +3.05%  '
$c.Style = "Normal"
